# "checking kinds of input files"
# Adds a new phoneme row (glottal stop "ʔ", feature-column label "Q")
# to the bottom of the Sheet2 feature-matrix table, and leaves the
# sheet scrolled/selected the way it was after typing that row in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Add the shared strings in the same order the author must have typed
# them (column B - the feature-label "Q" - before column A - the IPA
# glyph "ʔ") so they land in the shared-string table in that order.
$ws.Range("B43").Value = "Q"
$ws.Range("A43").Value = "ʔ"

# Feature values for the new row (glottal stop: consonantal + constrictedGlottis).
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("O43").Value = 0
$ws.Range("P43").Value = 1
$ws.Range("Q43").Value = 0
$ws.Range("R43").Value = 0
$ws.Range("S43").Value = 0
$ws.Range("T43").Value = 0
$ws.Range("X43").Value = 0
$ws.Range("Y43").Value = 0

# Restore the split-pane scroll position (top-left of the bottom-right
# pane moves from C2 to M16 once the view has scrolled down/right to
# the newly entered row).
$win = $excel.ActiveWindow
$win.SplitColumn = 12
$win.SplitRow = 15

# Final selection left active in the bottom-right pane after data entry.
$ws.Range("T43").Select()
